# "Generate Report for handoff"
#
# The localization-status report is regenerated: the source file's status
# flips from "Handoff transform failed" to "Ready for handoff", and a real
# handoff (.xlf) has now gone out for both the zh-cn and de-de targets, so
# each language sheet gains a "Latest Handoff File" hyperlink + timestamp
# and its "Handoff Reason" moves from "Ignored" to "Include".

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTest/oltest/blob/e6dd9c9c4069ad55807f26b3fec2c34f20c368b8"

# ---- Overview sheet: refresh the rolled-up status for the source file ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"

# ---- zh-cn sheet: record the new handoff file + timestamp, flip reason ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B2").Value = "Ready for handoff"

$zhcnHandoffFile = "5dfff4c6-4f74-471d-bf75-10703bd9bb67.07181534b274db6eb6e6948afdb588cc94d7efda.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("C2"), "$repoBase/e2e/$zhcnHandoffFile", "", "", $zhcnHandoffFile)

$zhcn.Range("D2").Value = "2016-01-20 03:49:40"
$zhcn.Range("H2").Value = "Include"

# ---- de-de sheet: same treatment ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B2").Value = "Ready for handoff"

$dedeHandoffFile = "5dfff4c6-4f74-471d-bf75-10703bd9bb67.07181534b274db6eb6e6948afdb588cc94d7efda.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("C2"), "$repoBase/e2e/$dedeHandoffFile", "", "", $dedeHandoffFile)

$dede.Range("D2").Value = "2016-01-20 03:49:51"
$dede.Range("H2").Value = "Include"
